$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.348.99"
$ws.Range("E2").Value = "  -6.48%  "
$ws.Range("D3").Value = "2.894.84"
$ws.Range("E3").Value = "  -4.16%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "550.42"
$ws.Range("E5").Value = "  -2.71%  "
$ws.Range("D6").Value = "123.28"
$ws.Range("E6").Value = "  -4.93%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").Value = "2.894.46"
$ws.Range("E8").Value = "  -4.23%  "
$ws.Range("D9").Value = "0.497"
$ws.Range("E9").Value = "  -0.40%  "
$ws.Range("E10").Value = "  -8.82%  "
$ws.Range("E11").Value = "  -10.94%  "
$ws.Range("D12").Value = "0.439"
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("E13").Value = "  -5.96%  "
$ws.Range("D14").Value = "32.51"
$ws.Range("E14").Value = "  -2.10%  "
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").Value = "3.371.76"
$ws.Range("E16").Value = "  -3.92%  "
$ws.Range("D17").Value = "2.894.63"
$ws.Range("E17").Value = "  -3.81%  "
$ws.Range("E18").Value = "  +4.17%  "
$ws.Range("D19").Value = "57.347.10"
$ws.Range("E19").Value = "  -6.40%  "
$ws.Range("D20").Value = "404.19"
$ws.Range("E20").Value = "  -8.31%  "
$ws.Range("D21").Value = "12.90"
$ws.Range("E21").Value = "  -2.70%  "
$ws.Range("D22").Value = "0.670"
$ws.Range("E22").Value = "  +0.45%  "
$ws.Range("D23").Value = "6.84"
$ws.Range("E23").Value = "  -4.87%  "
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").Value = "77.08"
$ws.Range("E25").Value = "  -2.85%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("E28").Value = "  -2.00%  "
$ws.Range("E29").Value = "  +1.66%  "
$ws.Range("E30").Value = "  -1.58%  "
$ws.Range("D31").Value = "6.02"
$ws.Range("E31").Value = "  -3.18%  "
$ws.Range("D32").Value = "24.67"
$ws.Range("E32").Value = "  -4.05%  "
$ws.Range("D33").Value = "0.0988"
$ws.Range("E33").Value = "  +4.39%  "
$ws.Range("E34").Value = "  -3.26%  "
$ws.Range("D35").Value = "0.906"
$ws.Range("E35").Value = "  -5.83%  "
$ws.Range("E36").Value = "  -12.88%  "
$ws.Range("D37").Value = "47.89"
$ws.Range("E37").Value = "  -5.15%  "
$ws.Range("D38").Value = "8.32"
$ws.Range("E38").Value = "  +6.54%  "
$ws.Range("E39").Value = "  -8.42%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.0341"
$ws.Range("E40").Value = "  -5.94%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "0.106"
$ws.Range("E41").Value = "  -3.21%  "
$ws.Range("D42").Value = "2.620.32"
$ws.Range("E42").Value = "  -2.35%  "
$ws.Range("D43").Value = "2.41"
$ws.Range("E43").Value = "  -3.35%  "
$ws.Range("D44").Value = "359.69"
$ws.Range("E44").Value = "  -5.43%  "
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "118.49"
$ws.Range("E46").Value = "  -0.54%  "
$ws.Range("E47").Value = "  -3.88%  "
$ws.Range("D48").Value = "0.107"
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("D49").Value = "1.94"
$ws.Range("E49").Value = "  -2.80%  "
$ws.Range("D50").Value = "22.84"
$ws.Range("E50").Value = "  -3.16%  "
$ws.Range("D51").Value = "1.95"
$ws.Range("E51").Value = "  -4.41%  "
